# Applies the diff: Tablas SubNetting.xlsx
# - T2 (sheet2): insert two new device rows (R1 F0/0 / R1 F1/0) before the
#   existing R1-F3/0 row, renumber PC1..PC6 interfaces to E0, extend table
#   "Tabla36" from A4:E11 to A4:E13.
# - T3 (sheet3): fix the device section to the 192.168.75.0/29 network
#   (PC7..PC10), add a trailing styled blank cell at C14.
# - Selections updated on all three sheets; T3 remains the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# T2 ("sheet2.xml")
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("T2")

# Row 5: F0/0 sub-interface of R1 (was: R1 / F3/0 / 192.168.47.1 / .240/28 / N/A)
$ws2.Range("A5").Value2 = "R1"
$ws2.Range("B5").Value2 = "F0/0"
$ws2.Range("C5").Value2 = "10.7.0.2"
$ws2.Range("D5").Value2 = "255.255.255.252 ó /30"
$ws2.Range("E5").Value2 = "N/A"

# Row 6: F1/0 sub-interface of R1
$ws2.Range("A6").Value2 = "R1"
$ws2.Range("B6").Value2 = "F1/0"
$ws2.Range("C6").Value2 = "10.7.0.6"
$ws2.Range("D6").Value2 = "255.255.255.252 ó /30"
$ws2.Range("E6").Value2 = "N/A"

# Row 7: the original R1/F3/0 row, shifted down
$ws2.Range("A7").Value2 = "R1"
$ws2.Range("B7").Value2 = "F3/0"
$ws2.Range("C7").Value2 = "192.168.47.1"
$ws2.Range("D7").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E7").Value2 = "N/A"

# Row 8: PC1, interface renumbered E1 -> E0
$ws2.Range("A8").Value2 = "PC1"
$ws2.Range("B8").Value2 = "E0"
$ws2.Range("C8").Value2 = "192.168.47.2"
$ws2.Range("D8").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E8").Value2 = "192.168.47.1"

# Row 9: PC2, interface renumbered E2 -> E0
$ws2.Range("A9").Value2 = "PC2"
$ws2.Range("B9").Value2 = "E0"
$ws2.Range("C9").Value2 = "192.168.47.3"
$ws2.Range("D9").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E9").Value2 = "192.168.47.1"

# Row 10: PC3, interface renumbered E3 -> E0
$ws2.Range("A10").Value2 = "PC3"
$ws2.Range("B10").Value2 = "E0"
$ws2.Range("C10").Value2 = "192.168.47.4"
$ws2.Range("D10").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E10").Value2 = "192.168.47.1"

# Row 11: PC4, interface renumbered E4 -> E0
$ws2.Range("A11").Value2 = "PC4"
$ws2.Range("B11").Value2 = "E0"
$ws2.Range("C11").Value2 = "192.168.47.5"
$ws2.Range("D11").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E11").Value2 = "192.168.47.1"

# Row 12 (new): PC5, E0
$ws2.Range("A12").Value2 = "PC5"
$ws2.Range("B12").Value2 = "E0"
$ws2.Range("C12").Value2 = "192.168.47.6"
$ws2.Range("D12").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E12").Value2 = "192.168.47.1"

# Row 13 (new): PC6, E0
$ws2.Range("A13").Value2 = "PC6"
$ws2.Range("B13").Value2 = "E0"
$ws2.Range("C13").Value2 = "192.168.47.7"
$ws2.Range("D13").Value2 = "255.255.255.240 ó /28"
$ws2.Range("E13").Value2 = "192.168.47.1"

# Grow the device table to cover the two new rows.
$lo2 = $ws2.ListObjects.Item("Tabla36")
$lo2.Resize($ws2.Range("A4:E13"))

# ---------------------------------------------------------------------
# T3 ("sheet3.xml")
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("T3")

# Row 5: F0/0 sub-interface of R1
$ws3.Range("A5").Value2 = "R1"
$ws3.Range("B5").Value2 = "F0/0"
$ws3.Range("C5").Value2 = "10.7.0.22"
$ws3.Range("D5").Value2 = "255.255.255.252 ó /30"
$ws3.Range("E5").Value2 = "N/A"

# Row 6: F1/0 sub-interface of R1
$ws3.Range("A6").Value2 = "R1"
$ws3.Range("B6").Value2 = "F1/0"
$ws3.Range("C6").Value2 = "10.7.0.18"
$ws3.Range("D6").Value2 = "255.255.255.252 ó /30"
$ws3.Range("E6").Value2 = "N/A"

# Row 7: F3/0 sub-interface of R1, now on the 192.168.75.0/29 network
$ws3.Range("A7").Value2 = "R1"
$ws3.Range("B7").Value2 = "F3/0"
$ws3.Range("C7").Value2 = "192.168.75.1"
$ws3.Range("D7").Value2 = "255.255.255.248 ó /29"
$ws3.Range("E7").Value2 = "N/A"

# Row 8: PC7
$ws3.Range("A8").Value2 = "PC7"
$ws3.Range("B8").Value2 = "E0"
$ws3.Range("C8").Value2 = "192.168.75.2"
$ws3.Range("D8").Value2 = "255.255.255.248 ó /29"
$ws3.Range("E8").Value2 = "192.168.75.1"

# Row 9: PC8
$ws3.Range("A9").Value2 = "PC8"
$ws3.Range("B9").Value2 = "E0"
$ws3.Range("C9").Value2 = "192.168.75.3"
$ws3.Range("D9").Value2 = "255.255.255.248 ó /29"
$ws3.Range("E9").Value2 = "192.168.75.1"

# Row 10: PC9
$ws3.Range("A10").Value2 = "PC9"
$ws3.Range("B10").Value2 = "E0"
$ws3.Range("C10").Value2 = "192.168.75.4"
$ws3.Range("D10").Value2 = "255.255.255.248 ó /29"
$ws3.Range("E10").Value2 = "192.168.75.1"

# Row 11: PC10
$ws3.Range("A11").Value2 = "PC10"
$ws3.Range("B11").Value2 = "E0"
$ws3.Range("C11").Value2 = "192.168.75.5"
$ws3.Range("D11").Value2 = "255.255.255.248 ó /29"
$ws3.Range("E11").Value2 = "192.168.75.1"

# Row 2 header-summary cells were pointing at shared strings that moved when
# E1..E5 were dropped from the shared-string table; re-assert the same text.
$ws3.Range("A2").Value2 = "192.168.75.0"
$ws3.Range("C2").Value2 = "192.168.75.1"
$ws3.Range("D2").Value2 = "192.168.75.6"
$ws3.Range("E2").Value2 = "192.168.75.7"

# New trailing row with a single styled (underlined, centered) blank cell.
$c14 = $ws3.Range("C14")
$c14.Value2 = ""
$c14.Font.Name = "Roboto Slab Light"
$c14.Font.Size = 12
$c14.Font.Underline = $true
$c14.HorizontalAlignment = -4108
$c14.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Selections (must be applied in this order so T3 ends up the active tab)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("T1")
$ws1.Range("A5").Select()

$ws2.Range("D7").Select()

$ws3.Range("B18").Select()
